$d = $word.ActiveDocument

# Locate the paragraph that holds the "LOB1036: ..." requisito text; the
# four paragraphs that follow it (an empty spacer paragraph, the
# "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph, another empty
# spacer paragraph, and an empty page-break paragraph) are being removed.
$marker = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "LOB1036:*") {
        $marker = $i
        break
    }
}

$startPara = $marker + 1
$endPara = $marker + 4

$rangeStart = $d.Paragraphs.Item($startPara).Range.Start
$rangeEnd = $d.Paragraphs.Item($endPara).Range.End

$r = $d.Range($rangeStart, $rangeEnd)
$r.Delete()
